# Update "想去人数" (column F) figures across all 4 sheets to match the
# refreshed scrape (gh-pages data regenerated at commit 456a3b4).
# Values below were derived by diffing old vs new F-column numbers per
# sheet/row against the workbook's current contents.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 744
$ws.Range("F3").Value = 14095
$ws.Range("F4").Value = 14095
$ws.Range("F5").Value = 14138
$ws.Range("F7").Value = 1390
$ws.Range("F8").Value = 5842
$ws.Range("F9").Value = 977
$ws.Range("F14").Value = 1527
$ws.Range("F15").Value = 430
$ws.Range("F17").Value = 1184
$ws.Range("F18").Value = 1800
$ws.Range("F20").Value = 33
$ws.Range("F21").Value = 2261
$ws.Range("F22").Value = 556
$ws.Range("F23").Value = 798
$ws.Range("F24").Value = 3283
$ws.Range("F26").Value = 306
$ws.Range("F27").Value = 2352
$ws.Range("F28").Value = 579
$ws.Range("F31").Value = 1772
$ws.Range("F32").Value = 1067
$ws.Range("F33").Value = 1358
$ws.Range("F34").Value = 98
$ws.Range("F35").Value = 138
$ws.Range("F36").Value = 4716
$ws.Range("F37").Value = 4781
$ws.Range("F38").Value = 297
$ws.Range("F40").Value = 668
$ws.Range("F41").Value = 678
$ws.Range("F42").Value = 3274
$ws.Range("F46").Value = 97
$ws.Range("F47").Value = 66
$ws.Range("F48").Value = 4414
$ws.Range("F49").Value = 560
$ws.Range("F50").Value = 280

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 45
$ws.Range("F4").Value = 116
$ws.Range("F15").Value = 19

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 7478
$ws.Range("F3").Value = 223
$ws.Range("F4").Value = 713

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 7478
$ws.Range("F3").Value = 744
$ws.Range("F4").Value = 223
$ws.Range("F5").Value = 713
$ws.Range("F6").Value = 14095
$ws.Range("F7").Value = 14138
$ws.Range("F9").Value = 1390
$ws.Range("F10").Value = 5842
$ws.Range("F11").Value = 977
$ws.Range("F12").Value = 116
$ws.Range("F15").Value = 1527
$ws.Range("F16").Value = 430
$ws.Range("F17").Value = 1184
$ws.Range("F18").Value = 1800
$ws.Range("F19").Value = 911
$ws.Range("F20").Value = 33
$ws.Range("F21").Value = 3283
$ws.Range("F22").Value = 306
$ws.Range("F23").Value = 2352
$ws.Range("F24").Value = 579
$ws.Range("F27").Value = 1772
$ws.Range("F29").Value = 19
$ws.Range("F31").Value = 1067
$ws.Range("F32").Value = 1358
$ws.Range("F33").Value = 98
$ws.Range("F34").Value = 4716
$ws.Range("F35").Value = 4781
$ws.Range("F36").Value = 297
$ws.Range("F38").Value = 668
$ws.Range("F39").Value = 678
$ws.Range("F40").Value = 3274
$ws.Range("F43").Value = 97
$ws.Range("F45").Value = 66
$ws.Range("F46").Value = 4414
$ws.Range("F47").Value = 560
$ws.Range("F48").Value = 280
